# ZBP_10_kontakty_aktivity.xlsx -- add the 12. 10. 2021 wave of data
# New column AJ on sheet "data" and new column AI on sheet "pocetR",
# plus the "aktualizace" date in each sheet's trailing caption row.

$wb  = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": append column AJ ("12. 10. 2021") after the last
# existing wave column AI.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# Header cell AJ1: copy AI1's formatting (bold/centered/bordered header
# style) onto AJ1, then set its value to the new wave's date label.
$ws1.Range("AI1").Copy()
$ws1.Range("AJ1").PasteSpecial(-4122)
$ws1.Range("AJ1").Value = "12. 10. 2021"

# Data rows 2-111: one % value per row for the new wave.
$data1 = @(0.17, 0.14, 0.18, 0.15, 0.36, 0.1, 0.07000000000000001, 0.18, 0.13, 0.52, 0.11, 0.08, 0.17, 0.19, 0.45, 0.26, 0.22, 0.18, 0.14, 0.2, 0.13, 0.06, 0.2, 0.09, 0.52, 0.2, 0.12, 0.18, 0.15, 0.35, 0.17, 0.24, 0.14, 0.23, 0.22, 0.12, 0.15, 0.16, 0.17, 0.4, 0.11, 0.16, 0.19, 0.22, 0.32, 0.2, 0.13, 0.17, 0.13, 0.37, 0.15, 0.12, 0.18, 0.16, 0.39, 0.2, 0.15, 0.17, 0.15, 0.33, 0.17, 0.13, 0.19, 0.12, 0.39, 0.19, 0.12, 0.19, 0.18, 0.32, 0.15, 0.15, 0.14, 0.2, 0.36, 0.15, 0.16, 0.14, 0.18, 0.37, 0.09, 0.05, 0.2, 0.13, 0.53, 0.09, 0.08, 0.15, 0.2, 0.48, 0.24, 0.21, 0.18, 0.14, 0.23, 0.12, 0.09, 0.17, 0.13, 0.49, 0.13, 0.09, 0.18, 0.18, 0.42, 0.28, 0.23, 0.17, 0.13, 0.19)
for ($i = 0; $i -lt $data1.Count; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 36).Value = $data1[$i]
}

# Trailing caption row (A112) references the update date - bump it.
$ws1.Range("A112").Value = "Život během pandemie, Kontakty vs. protektivní aktivity, % respondentů celkově a ve skupinách, aktualizace 20. 10. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": append column AI ("12. 10. 2021") after the last
# existing wave column AH.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AI1: copy AH1's formatting onto AI1, then set its value.
$ws2.Range("AH1").Copy()
$ws2.Range("AI1").PasteSpecial(-4122)
$ws2.Range("AI1").Value = "12. 10. 2021"

# Data rows 2-23: respondent counts for the new wave.
$data2 = @(1574, 362, 554, 658, 460, 681, 433, 253, 277, 1044, 790, 784, 803, 363, 192, 216, 188, 300, 303, 175, 254, 354)
for ($i = 0; $i -lt $data2.Count; $i++) {
    $row = 2 + $i
    $ws2.Cells.Item($row, 35).Value = $data2[$i]
}

# Row 24 is the blank trailer under the data block; keep column AI
# blank there too (mirrors the rest of that row).
$ws2.Range("AI24").Value = ""

# Trailing caption row (A24) references the update date - bump it.
$ws2.Range("A24").Value = "Život během pandemie, Kontakty vs. protektivní aktivity, velikost dotázaného souboru celkově a ve skupinách, aktualizace 20. 10. 2021"

Write-Output "edit complete"
